$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test-Cases")

# The "Approved/Rejected" column (I) is being marked "Approved" for every
# test case row. Format the data rows as Text (numFmtId 49, "@") first,
# then the header cell, matching the style order Excel records when the
# column body is formatted before its header.
$ws.Range("I2:I12").NumberFormat = "@"
$ws.Range("I1").NumberFormat = "@"

# Fill "Approved" down the whole Approved/Rejected data column.
$ws.Range("I2:I12").Value = "Approved"

# Leave the freshly-filled column selected, matching the saved selection.
$ws.Range("I1:I12").Select()
